$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refStyle = $ws.Range("D48").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.925.46'
$ws.Range("D2").Style = $refStyle
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.120.09'
$ws.Range("D3").Style = $refStyle
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.97'
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.39'
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = '  +1.84%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.112.59'
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("E10").Value = '  +9.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("E13").Value = '  +2.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.33'
$ws.Range("D14").Style = $refStyle
$ws.Range("E14").Value = '  +4.84%  '
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.636.92'
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.827.38'
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.116.42'
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '463.47'
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = '  +2.32%  '
$ws.Range("E21").Value = '  +1.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.729'
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = '  -0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("E24").Value = '  -3.20%  '
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  +8.54%  '
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.97'
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").Value = '  -2.07%  '
$ws.Range("E34").Value = '  +7.17%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +0.94%  '
$ws.Range("E37").Value = '  +9.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.06'
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.94'
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '446.23'
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = '  +4.87%  '
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.869.75'
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.275'
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  -0.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '35.77'
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = '  +3.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.29'
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = '  -1.12%  '
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("E51").Value = '  -1.20%  '
